$wb = $excel.ActiveWorkbook

# Rename the original (and only) sheet from "Product to search" to "Sheet"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sheet"

# Add a new row of data: header "expected card value on icon" with value 1
$ws1.Range("A3").Value = "expected card value on icon"
$ws1.Range("B3").Value = 1

# Add a new, empty worksheet named "Sheet1" right after the first sheet
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

# Make "Sheet" the active sheet again and select B3 as the active cell
$ws1.Activate()
$ws1.Range("B3").Select()
